$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 87: 2020-02-27
$ws.Range("A87").Value = 1582761600
$ws.Range("B87").NumberFormat = "@"
$ws.Range("B87").Value = "2020-02-27"
$ws.Range("B87").ClearFormats()
$ws.Range("C87").NumberFormat = "@"
$ws.Range("C87").Value = "5293"
$ws.Range("C87").ClearFormats()
$ws.Range("D87").Value = "AME"
$ws.Range("E87").Value = 1.8
$ws.Range("F87").Value = 1.81
$ws.Range("G87").Value = 1.73
$ws.Range("H87").Value = 1.75
$ws.Range("I87").Value = 1096200

# Row 88: 2020-02-28
$ws.Range("A88").Value = 1582848000
$ws.Range("B88").NumberFormat = "@"
$ws.Range("B88").Value = "2020-02-28"
$ws.Range("B88").ClearFormats()
$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "5293"
$ws.Range("C88").ClearFormats()
$ws.Range("D88").Value = "AME"
$ws.Range("E88").Value = 1.73
$ws.Range("F88").Value = 1.78
$ws.Range("G88").Value = 1.63
$ws.Range("H88").Value = 1.65
$ws.Range("I88").Value = 1103900
